# Sprint Backlog Template update:
#  - Remove the "DONE" column (G) marker text/values; header cell stays but empty.
#  - Append a new "Sprint 1 Report" block of backlog items (rows 59-75).
#  - Leave the final selection on the first new row's story cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Clear the old "DONE" tracking column ---
# Header cell (G1) keeps its style/position but no longer holds the "DONE" label.
$ws.Range("G1").ClearContents()
# All of the per-row DONE flags (G2:G58) are removed entirely.
$ws.Range("G2:G58").ClearContents()

# --- Append the new Sprint backlog rows (Unique ID, User story, Story Points, Assignee) ---
$newRows = @(
    @(6.7, "A Member should be able to add skills to his profile", 2, "Karim"),
    @(6.8, "A Member Should be able to add certificates to his profile", 2, "Moataz"),
    @(6.9, "A Member Should  be able to add completed tasks to his profile", 2, "Youssef"),
    @(7,   "A Member Should be able to add interests to his profile", 2, "Kashlan"),
    @(7.1, "A Member Should be able to add Past Events to his Profile", 2, "Mahmoud"),
    @(7.2, "A Partner Should be able to add Board Members to his profile", 2, "Basem"),
    @(7.3, "A Partner Should be able to add Partners to his profile", 2, "Sohail"),
    @(7.4, "A Partneer Should be able to add events to his profile", 2, "Hossam"),
    @(7.5, "A Partner Should be able to add project to his past projects", 2, "Gaafar"),
    @(7.6, "A Consultant Should be able Board Members to his Profile", 2, "Karim"),
    @(7.7, "A Consultant Should be able to add Partners to his profile", 2, "Youssef"),
    @(7.8, "A Consultant Should be able to add reports to his profile", 2, "Moataz"),
    @(7.9, "A Consultant Should be able to add events to his Profile", 2, "Kashlan"),
    @(8,   "An Educational Organization Should be able to add Courses to their Profile", 2, "Mahmoud"),
    @(8.1, "An Educational Organization Should be able to add Trainers to their Profile", 2, "Hossam"),
    @(8.2, "An Educational Organization Should be able to add Certificates to their Profile", 2, "Gaafar"),
    @(8.3, "An Educational Organization Should be able to add Training Programs to their Profile", 2, "Sohail")
)

$startRow = 59
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# --- Update the visible selection to the first new row (matches the re-saved view state) ---
$ws.Range("B59").Select()
